# Apply cryptos.xlsx symbol-list update (GitHub Actions scheduled refresh).
# Each entry forces the target cell to Text format before writing the new
# value so numeric-/percent-looking strings ("326.44", "-0.01%", ...) are
# preserved as literal text instead of being auto-parsed into numbers,
# matching the workbook's existing inline-string (text) cell storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '326.44' },
    @{ Cell = 'E2'; Value = '-0.01%' },
    @{ Cell = 'D3'; Value = '43.88' },
    @{ Cell = 'E3'; Value = '-1.80%' },
    @{ Cell = 'D4'; Value = '5.486' },
    @{ Cell = 'E4'; Value = '-1.13%' },
    @{ Cell = 'D5'; Value = '0.08015' },
    @{ Cell = 'E5'; Value = '-0.83%' },
    @{ Cell = 'D6'; Value = '1.980' },
    @{ Cell = 'E6'; Value = '4.28%' },
    @{ Cell = 'B7'; Value = 'GateToken' },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'D7'; Value = '4.291' },
    @{ Cell = 'E7'; Value = '-1.05%' },
    @{ Cell = 'B8'; Value = 'BTSEToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' },
    @{ Cell = 'D8'; Value = '2.568' },
    @{ Cell = 'E8'; Value = '-5.13%' },
    @{ Cell = 'B9'; Value = 'MXToken' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D9'; Value = '0.9530' },
    @{ Cell = 'E9'; Value = '0.45%' },
    @{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D10'; Value = '0.1122' },
    @{ Cell = 'E10'; Value = '-4.85%' },
    @{ Cell = 'B11'; Value = 'WazirX' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'D11'; Value = '0.1858' },
    @{ Cell = 'E11'; Value = '-1.85%' },
    @{ Cell = 'B12'; Value = 'MCDex' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' },
    @{ Cell = 'D12'; Value = '10.66' },
    @{ Cell = 'E12'; Value = '23.42%' },
    @{ Cell = 'B13'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D13'; Value = '0.09831' },
    @{ Cell = 'E13'; Value = '-2.79%' },
    @{ Cell = 'B14'; Value = 'BitrueCoin' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'D14'; Value = '0.04581' },
    @{ Cell = 'E14'; Value = '9.64%' },
    @{ Cell = 'B15'; Value = 'BitMartToken' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'D15'; Value = '0.1066' },
    @{ Cell = 'E15'; Value = '-0.01%' },
    @{ Cell = 'B16'; Value = 'BitForexToken' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'D16'; Value = '0.001260' },
    @{ Cell = 'E16'; Value = '-0.61%' },
    @{ Cell = 'B17'; Value = 'CoinExToken' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' },
    @{ Cell = 'D17'; Value = '0.04087' },
    @{ Cell = 'E17'; Value = '-3.82%' },
    @{ Cell = 'B18'; Value = 'TigerCash' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'D18'; Value = '0.005828' },
    @{ Cell = 'E18'; Value = '-4.50%' },
    @{ Cell = 'B19'; Value = 'LEO' },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'D19'; Value = '3.353' },
    @{ Cell = 'E19'; Value = '-6.92%' },
    @{ Cell = 'E20'; Value = '-0.29%' },
    @{ Cell = 'E21'; Value = '2.49%' },
    @{ Cell = 'D22'; Value = '0.2542' },
    @{ Cell = 'E22'; Value = '-4.55%' },
    @{ Cell = 'D23'; Value = '0.001251' },
    @{ Cell = 'E23'; Value = '1.09%' },
    @{ Cell = 'D24'; Value = '0.004332' },
    @{ Cell = 'E24'; Value = '-5.94%' },
    @{ Cell = 'D25'; Value = '0.0001161' },
    @{ Cell = 'E25'; Value = '-5.74%' },
    @{ Cell = 'D26'; Value = '0.0003738' },
    @{ Cell = 'E26'; Value = '-6.53%' },
    @{ Cell = 'D38'; Value = '0.02562' },
    @{ Cell = 'E38'; Value = '-3.71%' },
    @{ Cell = 'D39'; Value = '0.05679' },
    @{ Cell = 'E39'; Value = '2.17%' },
    @{ Cell = 'D40'; Value = '0.007550' },
    @{ Cell = 'E40'; Value = '-1.83%' },
    @{ Cell = 'D41'; Value = '0.1395' },
    @{ Cell = 'E41'; Value = '0.02%' },
    @{ Cell = 'D42'; Value = '0.007603' },
    @{ Cell = 'E42'; Value = '-32.92%' },
    @{ Cell = 'D43'; Value = '0.002011' },
    @{ Cell = 'E43'; Value = '-2.26%' },
    @{ Cell = 'D44'; Value = '0.008504' },
    @{ Cell = 'E44'; Value = '-7.53%' },
    @{ Cell = 'E45'; Value = '-0.32%' },
    @{ Cell = 'E46'; Value = '-0.41%' },
    @{ Cell = 'E47'; Value = '55.09%' },
    @{ Cell = 'D48'; Value = '0.003086' },
    @{ Cell = 'E48'; Value = '-10.22%' },
    @{ Cell = 'D49'; Value = '0.00002096' },
    @{ Cell = 'E49'; Value = '-0.41%' },
    @{ Cell = 'D50'; Value = '0.0001997' },
    @{ Cell = 'E50'; Value = '-0.41%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
